# adding files to server pressing changes to date_log
# Populate the "date_logger" sheet with the employee hours/day log and make
# it the active sheet (matches the committed diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # date_logger

# --- Row data: emp_id (B), date (C), hours_worked (D) ---------------------
$empIds = @(1001, 1002, 1003, 1004, 1005)
$dates  = @(42467, 42468, 42471, 42472)   # 2016-04-07/08/11/12 (Thu/Fri/Mon/Tue)

$row = 2
foreach ($d in $dates) {
    foreach ($e in $empIds) {
        $ws.Cells.Item($row, 2).Value = $e
        $ws.Cells.Item($row, 3).Value = $d
        $row = $row + 1
    }
}

$ws.Range("D2").Value = 9
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 7
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 9

$ws.Range("D7").Value = 6
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 8
$ws.Range("D10").Value = 9
$ws.Range("D11").Value = 11

$ws.Range("D12").Value = 4
$ws.Range("D13").Value = 7
$ws.Range("D14").Value = 8
$ws.Range("D15").Value = 0
$ws.Range("D16").Value = 10

$ws.Range("D17").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("D19").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("D21").Value = 9

# --- Day-of-week labels in column K, one per date block --------------------
$ws.Range("K1").Value = "day"
$ws.Range("K2").Value = "Thrs"
$ws.Range("K7").Value = "Fri"
$ws.Range("K12").Value = "Monday"
$ws.Range("K17").Value = "Tuesday"

# --- Date formatting: stamp one cell with the date number format, then -----
# --- propagate that exact style to the rest of column C (C2:C21) so all ----
# --- date cells share a single style record instead of duplicating one -----
# --- per cell. ---------------------------------------------------------
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Sheet view: make date_logger the active/selected tab with F19 selected
$ws.Activate()
$ws.Range("F19").Select() | Out-Null

# --- Page setup (paper size + orientation) ---------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
